# Generate Report for Handback
# The 96721e3b file has now also been handed back (was "Ready for handoff"),
# so the report is regenerated: row 2 now carries the 96721e3b identity and
# row 3 carries the af3db7a1 identity, both marked "Handed back: in sync
# with en-US", and fresh handback timestamps are recorded for the zh-cn /
# de-de detail sheets. Hyperlink targets stay anchored to their sheet
# position; only their displayed text is refreshed to match the new values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = "96721e3b-4bbb-400b-87de-c124531acbca.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "2016-30-17 03:30:26"

$ws.Range("A3").Value = "af3db7a1-60b9-45b8-b4be-746fe66c8af5.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "2016-30-17 03:30:26"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/70e9314ab376397f21f9fd1ce9313810b27abc39/e2e/af3db7a1-60b9-45b8-b4be-746fe66c8af5.md", "", "", "96721e3b-4bbb-400b-87de-c124531acbca.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d13dbbdc6d1ee0e37fce01d94483358e4d0f953e/e2e/96721e3b-4bbb-400b-87de-c124531acbca.md", "", "", "af3db7a1-60b9-45b8-b4be-746fe66c8af5.md") | Out-Null

$ws.Range("A2").Style = "HyperLink"
$ws.Range("A3").Style = "HyperLink"

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = "96721e3b-4bbb-400b-87de-c124531acbca.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "96721e3b-4bbb-400b-87de-c124531acbca.3fbb8e26b1e70ab87987783469d7741a664c0a75.zh-cn.xlf"
$ws.Range("E2").Value = "2016-03-17 03:30:17"
$ws.Range("F2").Value = "96721e3b-4bbb-400b-87de-c124531acbca.md"
$ws.Range("G2").Value = "96721e3b-4bbb-400b-87de-c124531acbca.3fbb8e26b1e70ab87987783469d7741a664c0a75.zh-cn.xlf"
$ws.Range("H2").Value = "2016-03-17 03:31:03"
$ws.Range("I2").Value = "Include"

$ws.Range("A3").Value = "af3db7a1-60b9-45b8-b4be-746fe66c8af5.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "af3db7a1-60b9-45b8-b4be-746fe66c8af5.f483e1e6dce6bb521126d5fbb3ba6eec7609aa7c.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-17 03:30:17"
$ws.Range("F3").Value = "af3db7a1-60b9-45b8-b4be-746fe66c8af5.md"
$ws.Range("G3").Value = "af3db7a1-60b9-45b8-b4be-746fe66c8af5.f483e1e6dce6bb521126d5fbb3ba6eec7609aa7c.zh-cn.xlf"
$ws.Range("H3").Value = "2016-03-17 03:31:21"
$ws.Range("I3").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/70e9314ab376397f21f9fd1ce9313810b27abc39/e2e/af3db7a1-60b9-45b8-b4be-746fe66c8af5.md", "", "", "96721e3b-4bbb-400b-87de-c124531acbca.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/70e9314ab376397f21f9fd1ce9313810b27abc39/e2e/af3db7a1-60b9-45b8-b4be-746fe66c8af5.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2aa3d597750e444e0c08c0527218279e725e28bf/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/mt/af3db7a1-60b9-45b8-b4be-746fe66c8af5.f483e1e6dce6bb521126d5fbb3ba6eec7609aa7c.zh-cn.xlf", "", "", "96721e3b-4bbb-400b-87de-c124531acbca.3fbb8e26b1e70ab87987783469d7741a664c0a75.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e36b10481ea8bbd5660a92acf3e64a3dd9eb36f0/e2e/af3db7a1-60b9-45b8-b4be-746fe66c8af5.md", "", "", "96721e3b-4bbb-400b-87de-c124531acbca.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/da34c2b70d416ed5d97015b67fb5ac8a5aafaee1/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/af3db7a1-60b9-45b8-b4be-746fe66c8af5.f483e1e6dce6bb521126d5fbb3ba6eec7609aa7c.zh-cn.xlf", "", "", "96721e3b-4bbb-400b-87de-c124531acbca.3fbb8e26b1e70ab87987783469d7741a664c0a75.zh-cn.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d13dbbdc6d1ee0e37fce01d94483358e4d0f953e/e2e/96721e3b-4bbb-400b-87de-c124531acbca.md", "", "", "af3db7a1-60b9-45b8-b4be-746fe66c8af5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/d13dbbdc6d1ee0e37fce01d94483358e4d0f953e/e2e/96721e3b-4bbb-400b-87de-c124531acbca.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2aa3d597750e444e0c08c0527218279e725e28bf/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/mt/96721e3b-4bbb-400b-87de-c124531acbca.3fbb8e26b1e70ab87987783469d7741a664c0a75.zh-cn.xlf", "", "", "af3db7a1-60b9-45b8-b4be-746fe66c8af5.f483e1e6dce6bb521126d5fbb3ba6eec7609aa7c.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e36b10481ea8bbd5660a92acf3e64a3dd9eb36f0/e2e/96721e3b-4bbb-400b-87de-c124531acbca.md", "", "", "af3db7a1-60b9-45b8-b4be-746fe66c8af5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/da34c2b70d416ed5d97015b67fb5ac8a5aafaee1/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/96721e3b-4bbb-400b-87de-c124531acbca.3fbb8e26b1e70ab87987783469d7741a664c0a75.zh-cn.xlf", "", "", "af3db7a1-60b9-45b8-b4be-746fe66c8af5.f483e1e6dce6bb521126d5fbb3ba6eec7609aa7c.zh-cn.xlf") | Out-Null

foreach ($addr in @("A2","B2","D2","F2","G2","A3","B3","D3","F3","G3")) {
    $ws.Range($addr).Style = "HyperLink"
}

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = "96721e3b-4bbb-400b-87de-c124531acbca.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "96721e3b-4bbb-400b-87de-c124531acbca.3fbb8e26b1e70ab87987783469d7741a664c0a75.de-de.xlf"
$ws.Range("E2").Value = "2016-03-17 03:30:26"
$ws.Range("F2").Value = "96721e3b-4bbb-400b-87de-c124531acbca.md"
$ws.Range("G2").Value = "96721e3b-4bbb-400b-87de-c124531acbca.3fbb8e26b1e70ab87987783469d7741a664c0a75.de-de.xlf"
$ws.Range("H2").Value = "2016-03-17 03:31:21"
$ws.Range("I2").Value = "Include"

$ws.Range("A3").Value = "af3db7a1-60b9-45b8-b4be-746fe66c8af5.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "af3db7a1-60b9-45b8-b4be-746fe66c8af5.f483e1e6dce6bb521126d5fbb3ba6eec7609aa7c.de-de.xlf"
$ws.Range("E3").Value = "2016-03-17 03:30:26"
$ws.Range("F3").Value = "af3db7a1-60b9-45b8-b4be-746fe66c8af5.md"
$ws.Range("G3").Value = "af3db7a1-60b9-45b8-b4be-746fe66c8af5.f483e1e6dce6bb521126d5fbb3ba6eec7609aa7c.de-de.xlf"
$ws.Range("H3").Value = "2016-03-17 03:29:36"
$ws.Range("I3").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/70e9314ab376397f21f9fd1ce9313810b27abc39/e2e/af3db7a1-60b9-45b8-b4be-746fe66c8af5.md", "", "", "96721e3b-4bbb-400b-87de-c124531acbca.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/70e9314ab376397f21f9fd1ce9313810b27abc39/e2e/af3db7a1-60b9-45b8-b4be-746fe66c8af5.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/20bbcd2ccfa670940b7ecb96b6afad752caca91e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/mt/af3db7a1-60b9-45b8-b4be-746fe66c8af5.f483e1e6dce6bb521126d5fbb3ba6eec7609aa7c.de-de.xlf", "", "", "96721e3b-4bbb-400b-87de-c124531acbca.3fbb8e26b1e70ab87987783469d7741a664c0a75.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/317e443c2f647e4d52fe323194125e597d7b9458/e2e/af3db7a1-60b9-45b8-b4be-746fe66c8af5.md", "", "", "96721e3b-4bbb-400b-87de-c124531acbca.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a57ecaee65058d55c5ed8e62bf2279bbdff2359b/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/af3db7a1-60b9-45b8-b4be-746fe66c8af5.f483e1e6dce6bb521126d5fbb3ba6eec7609aa7c.de-de.xlf", "", "", "96721e3b-4bbb-400b-87de-c124531acbca.3fbb8e26b1e70ab87987783469d7741a664c0a75.de-de.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d13dbbdc6d1ee0e37fce01d94483358e4d0f953e/e2e/96721e3b-4bbb-400b-87de-c124531acbca.md", "", "", "af3db7a1-60b9-45b8-b4be-746fe66c8af5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/d13dbbdc6d1ee0e37fce01d94483358e4d0f953e/e2e/96721e3b-4bbb-400b-87de-c124531acbca.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/20bbcd2ccfa670940b7ecb96b6afad752caca91e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/mt/96721e3b-4bbb-400b-87de-c124531acbca.3fbb8e26b1e70ab87987783469d7741a664c0a75.de-de.xlf", "", "", "af3db7a1-60b9-45b8-b4be-746fe66c8af5.f483e1e6dce6bb521126d5fbb3ba6eec7609aa7c.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/317e443c2f647e4d52fe323194125e597d7b9458/e2e/96721e3b-4bbb-400b-87de-c124531acbca.md", "", "", "af3db7a1-60b9-45b8-b4be-746fe66c8af5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a57ecaee65058d55c5ed8e62bf2279bbdff2359b/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/96721e3b-4bbb-400b-87de-c124531acbca.3fbb8e26b1e70ab87987783469d7741a664c0a75.de-de.xlf", "", "", "af3db7a1-60b9-45b8-b4be-746fe66c8af5.f483e1e6dce6bb521126d5fbb3ba6eec7609aa7c.de-de.xlf") | Out-Null

foreach ($addr in @("A2","B2","D2","F2","G2","A3","B3","D3","F3","G3")) {
    $ws.Range($addr).Style = "HyperLink"
}

$ws = $wb.Worksheets.Item("Overview")
$ws.Select()
$ws.Range("A1").Select()
